# Edit script: add new survey wave column (24. 8. 2021) to both sheets
# Sheet 1 "data": insert column AG with percentage values
# Sheet 2 "pocetR": insert column AF with sample-size counts
# Also refresh the "aktualizace" (updated on) date in the footer rows.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# --- Sheet 1 ("data") ---------------------------------------------------
# Header cell AG1: copy formatting from the previous header cell (AF1) so the
# new column matches the existing header style, then set its value.
$ws1.Range("AF1").Copy() | Out-Null
$ws1.Range("AG1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws1.Range("AG1").Value = "24. 8. 2021"

# Data rows 2-61: new percentage values for the 24. 8. 2021 wave.
$ws1.Range("AG2").Value = 0.5600000000000001
$ws1.Range("AG3").Value = 0.29
$ws1.Range("AG4").Value = 0.15
$ws1.Range("AG5").Value = 0.37
$ws1.Range("AG6").Value = 0.33
$ws1.Range("AG7").Value = 0.3
$ws1.Range("AG8").Value = 0.58
$ws1.Range("AG9").Value = 0.29
$ws1.Range("AG10").Value = 0.13
$ws1.Range("AG11").Value = 0.58
$ws1.Range("AG12").Value = 0.29
$ws1.Range("AG13").Value = 0.13
$ws1.Range("AG14").Value = 0.48
$ws1.Range("AG15").Value = 0.3
$ws1.Range("AG16").Value = 0.22
$ws1.Range("AG17").Value = 0.59
$ws1.Range("AG18").Value = 0.28
$ws1.Range("AG19").Value = 0.13
$ws1.Range("AG20").Value = 0.48
$ws1.Range("AG21").Value = 0.31
$ws1.Range("AG22").Value = 0.21
$ws1.Range("AG23").Value = 0.47
$ws1.Range("AG24").Value = 0.35
$ws1.Range("AG25").Value = 0.18
$ws1.Range("AG26").Value = 0.48
$ws1.Range("AG27").Value = 0.31
$ws1.Range("AG28").Value = 0.21
$ws1.Range("AG29").Value = 0.59
$ws1.Range("AG30").Value = 0.28
$ws1.Range("AG31").Value = 0.13
$ws1.Range("AG32").Value = 0.62
$ws1.Range("AG33").Value = 0.27
$ws1.Range("AG34").Value = 0.11
$ws1.Range("AG35").Value = 0.35
$ws1.Range("AG36").Value = 0.41
$ws1.Range("AG37").Value = 0.24
$ws1.Range("AG38").Value = 0.48
$ws1.Range("AG39").Value = 0.32
$ws1.Range("AG40").Value = 0.2
$ws1.Range("AG41").Value = 0.66
$ws1.Range("AG42").Value = 0.23
$ws1.Range("AG43").Value = 0.11
$ws1.Range("AG44").Value = 0.6899999999999999
$ws1.Range("AG45").Value = 0.22
$ws1.Range("AG46").Value = 0.09
$ws1.Range("AG47").Value = 0.55
$ws1.Range("AG48").Value = 0.3
$ws1.Range("AG49").Value = 0.15
$ws1.Range("AG50").Value = 0.68
$ws1.Range("AG51").Value = 0.27
$ws1.Range("AG52").Value = 0.05
$ws1.Range("AG53").Value = 0.49
$ws1.Range("AG54").Value = 0.33
$ws1.Range("AG55").Value = 0.18
$ws1.Range("AG56").Value = 0.6899999999999999
$ws1.Range("AG57").Value = 0.22
$ws1.Range("AG58").Value = 0.09
$ws1.Range("AG59").Value = 0.6899999999999999
$ws1.Range("AG60").Value = 0.18
$ws1.Range("AG61").Value = 0.13

# Footer row 62: bump the "aktualizace" date shown in the title text.
$ws1.Range("A62").Value = "Život během pandemie, Obavy ze ztráty práce, % respondentů celkově a ve skupinách, aktualizace 1. 9. 2021"

# --- Sheet 2 ("pocetR") ---------------------------------------------------
# Header cell AF1: copy formatting from the previous header cell (AE1) so the
# new column matches the existing header style, then set its value.
$ws2.Range("AE1").Copy() | Out-Null
$ws2.Range("AF1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
$ws2.Range("AF1").Value = "24. 8. 2021"

# Data rows 2-24: new sample-size counts for the 24. 8. 2021 wave.
$ws2.Range("AF2").Value = 1073
$ws2.Range("AF3").Value = 93
$ws2.Range("AF4").Value = 980
$ws2.Range("AF5").Value = 837
$ws2.Range("AF6").Value = 152
$ws2.Range("AF7").Value = 10
$ws2.Range("AF8").Value = 74
$ws2.Range("AF9").Value = 799
$ws2.Range("AF10").Value = 142
$ws2.Range("AF11").Value = 64
$ws2.Range("AF12").Value = 68
$ws2.Range("AF13").Value = 387
$ws2.Range("AF14").Value = 436
$ws2.Range("AF15").Value = 250
$ws2.Range("AF16").Value = 121
$ws2.Range("AF17").Value = 335
$ws2.Range("AF18").Value = 332
$ws2.Range("AF19").Value = 180
$ws2.Range("AF20").Value = 293
$ws2.Range("AF21").Value = 91
$ws2.Range("AF22").Value = 292
$ws2.Range("AF23").Value = 158
$ws2.Range("AF24").Value = 85

# Footer row 25: bump the "aktualizace" date, and extend the trailing blank
# cell range (AF25) to match the new last column, mirroring the empty
# formatting-only cells B25:AE25.
$ws2.Range("A25").Value = "Život během pandemie, Obavy ze ztráty práce, velikost dotázaného souboru celkově a ve skupinách, aktualizace 1. 9. 2021"
$ws2.Range("AE25").Copy() | Out-Null
$ws2.Range("AF25").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
